$d = $word.ActiveDocument

function Get-ParaByText {
    param([string]$text)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13) -eq $text) {
            return $p
        }
    }
    return $null
}

function Insert-BlankParaAfter {
    # Inserts an empty <w:p/> immediately after the given paragraph.
    param($para)
    $r = $para.Range
    $r.Collapse(0)   # wdCollapseEnd
    $r.InsertAfter([char]13)
}

function Insert-TextParaAfter {
    # Inserts a new paragraph with $text immediately after $para, returns nothing;
    # relies on the blank paragraph + fill-in-place pattern so no stray runs are left behind.
    param($para, [string]$text)
    Insert-BlankParaAfter $para
    $newPara = $para.Next()
    $newPara.Range.InsertBefore($text)
}

function Set-HeaderStyle {
    # Applies bold + 6pt before/after spacing to a paragraph's visible text
    # (excludes the paragraph mark so <w:rPr> stays only on the <w:r>, not on <w:pPr>).
    param($para)
    $para.SpaceBefore = 6
    $para.SpaceAfter = 6
    $r = $para.Range
    $r.MoveEnd(1, -1)  # wdCharacter, exclude trailing paragraph mark
    $r.Bold = 1
}

# ---------------------------------------------------------------------------
# 1) "2. Customer Research..." section: blank line before, bold+spacing header,
#    blank line after, then two new supporting paragraphs.
# ---------------------------------------------------------------------------
$pProblem = Get-ParaByText("[EM] Consolidation - Finalize all linked engagements at the same time (Epic 4233310)1. Customer Problemeach engagement individually is slow and increases risk of mismatch across tiers.")
Insert-BlankParaAfter $pProblem

$pResearch = Get-ParaByText("2. Customer Researchfeedback: desire for bulk finalize from parent with controls.")
Set-HeaderStyle $pResearch

Insert-BlankParaAfter $pResearch
$pResearch2 = Get-ParaByText("2. Customer Researchfeedback: desire for bulk finalize from parent with controls.")
Insert-TextParaAfter $pResearch2 "This capability was requested as feedback from an enterprise-level accounting firm, reflecting needs observed in large multi-entity audit workflows."

$pFirm = Get-ParaByText("This capability was requested as feedback from an enterprise-level accounting firm, reflecting needs observed in large multi-entity audit workflows.")
Insert-TextParaAfter $pFirm "We are also building this to achieve competitive parity with Wolters Kluwer ProSystem fx Engagement, which offers similar functionality."

# ---------------------------------------------------------------------------
# 2) "3. Our Solution..." header: blank before/after + bold+spacing.
# ---------------------------------------------------------------------------
$pParity = Get-ParaByText("We are also building this to achieve competitive parity with Wolters Kluwer ProSystem fx Engagement, which offers similar functionality.")
Insert-BlankParaAfter $pParity

$pSolution = Get-ParaByText("3. Our SolutionBulk Finalization wizard with scope selection and prerequisites checks.")
Set-HeaderStyle $pSolution
Insert-BlankParaAfter $pSolution

# ---------------------------------------------------------------------------
# 3) "4. Product Metrics..." header: blank before/after + bold+spacing.
# ---------------------------------------------------------------------------
$pDependencies = Get-ParaByText("Dependencies: finalization services, permissions.")
Insert-BlankParaAfter $pDependencies

$pMetrics = Get-ParaByText("4. Product MetricsTime to finalize 5 children: -60%.")
Set-HeaderStyle $pMetrics
Insert-BlankParaAfter $pMetrics

# ---------------------------------------------------------------------------
# 4) "Appendix: Quick prototype" heading: blank before/after + bold+spacing
#    (keeps its Heading2 pStyle).
# ---------------------------------------------------------------------------
$pLinks = Get-ParaByText("Appendix: LinksHYPERLINK `"https://dev.azure.com/tr-tax/TaxProf/_workitems/edit/4233310`"Open Epic 4233310")
Insert-BlankParaAfter $pLinks

$pPrototype = Get-ParaByText("Appendix: Quick prototype")
Set-HeaderStyle $pPrototype
Insert-BlankParaAfter $pPrototype

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
